$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the description text for each key in column A (rows 2-6)
$ws.Range("A2").Value = "The ""hint"" text used when prompting for the uer's login name"
$ws.Range("A3").Value = """Ok"" text displayed on dialog buttons"
$ws.Range("A4").Value = """Cancel"" text displayed on dialog buttons"
$ws.Range("A5").Value = """Yes"" text displayed on dialog buttons when making a request from the user"
$ws.Range("A6").Value = """No"" text displayed on dialog buttons when making a request from the user"

# Widen column A to fit the new descriptions
# (target stored width is 43.83203125 characters; this runtime snaps
# ColumnWidth to the nearest 1/6-character pixel grid, so 43.0 is the
# input that lands closest on the achievable grid -> 43.8333... stored)
$ws.Columns.Item(1).ColumnWidth = 43.0

# Update the active selection to A7
$ws.Range("A7").Select()
